$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5032.4116
$ws.Range("I98").Value = 2773.6155
$ws.Range("K98").Value = 2773.6155
$ws.Range("M98").Value = -1275.6155
$ws.Range("H122").Value = 5032.4116
$ws.Range("I122").Value = 2773.6155
$ws.Range("K122").Value = 8320.8465
$ws.Range("M122").Value = -5870.8465
$ws.Range("H125").Value = 2356.5
$ws.Range("J125").Value = 1963
$ws.Range("L125").Value = 17667
$ws.Range("N125").Value = -22587
$ws.Range("H138").Value = 2162.3416
$ws.Range("I138").Value = 1355.7
$ws.Range("K138").Value = 4067.1
$ws.Range("M138").Value = 1072.9
$ws.Range("H140").Value = 66175.21000000001
$ws.Range("I140").Value = 50000
$ws.Range("J140").Value = 68871.086
$ws.Range("K140").Value = 50000
$ws.Range("L140").Value = 68871.086
$ws.Range("M140").Value = -44820
$ws.Range("N140").Value = -79231.086

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2389.8713
$ws.Range("I32").Value = 1901.0938
$ws.Range("K32").Value = 1901.0938
$ws.Range("M32").Value = -1614.0938
$ws.Range("H45").Value = 3333
$ws.Range("I45").Value = 5999
$ws.Range("K45").Value = 5999
$ws.Range("M45").Value = -5622
$ws.Range("H61").Value = 7178.25
$ws.Range("I61").Value = 6499.5
$ws.Range("J61").Value = 7857
$ws.Range("K61").Value = 6499.5
$ws.Range("L61").Value = 7857
$ws.Range("M61").Value = -6287.5
$ws.Range("N61").Value = -8281
$ws.Range("H123").Value = 71500
$ws.Range("J123").Value = 71500
$ws.Range("L123").Value = 71500
$ws.Range("N123").Value = -81300
$ws.Range("H136").Value = 7178.25
$ws.Range("I136").Value = 6499.5
$ws.Range("J136").Value = 7857
$ws.Range("K136").Value = 19498.5
$ws.Range("L136").Value = 23571
$ws.Range("M136").Value = -16948.5
$ws.Range("N136").Value = -28671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2251.4
$ws.Range("J105").Value = 2997.5
$ws.Range("L105").Value = 2997.5
$ws.Range("N105").Value = -6491.5
$ws.Range("H107").Value = 1655.3846
$ws.Range("I107").Value = 1718.3334
$ws.Range("K107").Value = 1718.3334
$ws.Range("M107").Value = 201.6666
$ws.Range("H116").Value = 39750
$ws.Range("J116").Value = 39750
$ws.Range("L116").Value = 39750
$ws.Range("N116").Value = -48928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4549.263
$ws.Range("I31").Value = 1187.3334
$ws.Range("J31").Value = 5592.6206
$ws.Range("K31").Value = 1187.3334
$ws.Range("L31").Value = 5592.6206
$ws.Range("M31").Value = -892.3334
$ws.Range("N31").Value = -6182.6206
$ws.Range("H34").Value = 4549.263
$ws.Range("I34").Value = 1187.3334
$ws.Range("J34").Value = 5592.6206
$ws.Range("K34").Value = 1187.3334
$ws.Range("L34").Value = 5592.6206
$ws.Range("M34").Value = -985.3334
$ws.Range("N34").Value = -5996.6206
$ws.Range("H52").Value = 118657.8
$ws.Range("J52").Value = 118657.8
$ws.Range("L52").Value = 118657.8
$ws.Range("N52").Value = -119245.8
$ws.Range("H62").Value = 2455.4
$ws.Range("I62").Value = 2319.25
$ws.Range("K62").Value = 2319.25
$ws.Range("M62").Value = -1695.25
$ws.Range("H65").Value = 2455.4
$ws.Range("I65").Value = 2319.25
$ws.Range("K65").Value = 11596.25
$ws.Range("M65").Value = -8476.25
$ws.Range("H122").Value = 2319
$ws.Range("I122").Value = 2302.5454
$ws.Range("K122").Value = 6907.6362
$ws.Range("M122").Value = -4457.6362
$ws.Range("H134").Value = 2726.2666
$ws.Range("I134").Value = 2489.9
$ws.Range("K134").Value = 7469.700000000001
$ws.Range("M134").Value = -4934.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7001.2
$ws.Range("J75").Value = 7001.2
$ws.Range("L75").Value = 21003.6
$ws.Range("N75").Value = -22999.6
$ws.Range("H78").Value = 7001.2
$ws.Range("J78").Value = 7001.2
$ws.Range("L78").Value = 63010.8
$ws.Range("N78").Value = -72994.79999999999
$ws.Range("H98").Value = 631.5833
$ws.Range("J98").Value = 664.44446
$ws.Range("L98").Value = 1993.33338
$ws.Range("N98").Value = -4989.33338
$ws.Range("H107").Value = 1357.3077
$ws.Range("I107").Value = 1197.25
$ws.Range("J107").Value = 1613.4
$ws.Range("K107").Value = 3591.75
$ws.Range("L107").Value = 4840.200000000001
$ws.Range("M107").Value = -1671.75
$ws.Range("N107").Value = -8680.200000000001
$ws.Range("H122").Value = 935.4783
$ws.Range("I122").Value = 841.75
$ws.Range("K122").Value = 7575.75
$ws.Range("M122").Value = -5125.75
$ws.Range("H131").Value = 3126575.8
$ws.Range("J131").Value = 1962418.9
$ws.Range("L131").Value = 5887256.699999999
$ws.Range("N131").Value = -5897336.699999999
$ws.Range("H137").Value = 3023.9412
$ws.Range("I137").Value = 1488.5
$ws.Range("J137").Value = 4388.778
$ws.Range("K137").Value = 4465.5
$ws.Range("L137").Value = 13166.334
$ws.Range("M137").Value = 634.5
$ws.Range("N137").Value = -23366.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3316.7778
$ws.Range("I102").Value = 3097.6155
$ws.Range("K102").Value = 3097.6155
$ws.Range("M102").Value = -1475.6155
$ws.Range("H113").Value = 3527.0667
$ws.Range("I113").Value = 3504
$ws.Range("J113").Value = 3542.4443
$ws.Range("K113").Value = 3504
$ws.Range("L113").Value = 3542.4443
$ws.Range("M113").Value = -1334
$ws.Range("N113").Value = -7882.4443
$ws.Range("H122").Value = 7721.722
$ws.Range("I122").Value = 9142.286
$ws.Range("J122").Value = 2749.75
$ws.Range("K122").Value = 27426.858
$ws.Range("L122").Value = 8249.25
$ws.Range("M122").Value = -24976.858
$ws.Range("N122").Value = -13149.25
$ws.Range("H126").Value = 11477581
$ws.Range("I126").Value = 9481.166999999999
$ws.Range("J126").Value = 17732908
$ws.Range("K126").Value = 28443.501
$ws.Range("L126").Value = 53198724
$ws.Range("M126").Value = -25973.501
$ws.Range("N126").Value = -53203664
$ws.Range("H132").Value = 8397
$ws.Range("J132").Value = 8699.315000000001
$ws.Range("L132").Value = 26097.945
$ws.Range("N132").Value = -31157.945

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10361.0625
$ws.Range("I61").Value = 9290.615
$ws.Range("K61").Value = 9290.615
$ws.Range("M61").Value = -9088.615
$ws.Range("H68").Value = 3860.7
$ws.Range("I68").Value = 1481.5
$ws.Range("K68").Value = 1481.5
$ws.Range("M68").Value = -732.5
$ws.Range("H71").Value = 3860.7
$ws.Range("I71").Value = 1481.5
$ws.Range("K71").Value = 7407.5
$ws.Range("M71").Value = -3663.5
$ws.Range("H82").Value = 4438.25
$ws.Range("I82").Value = 1299.3334
$ws.Range("J82").Value = 6321.6
$ws.Range("K82").Value = 1299.3334
$ws.Range("L82").Value = 6321.6
$ws.Range("M82").Value = -938.3334
$ws.Range("N82").Value = -7043.6
$ws.Range("H85").Value = 4438.25
$ws.Range("I85").Value = 1299.3334
$ws.Range("J85").Value = 6321.6
$ws.Range("K85").Value = 1299.3334
$ws.Range("L85").Value = 6321.6
$ws.Range("M85").Value = -51.33339999999998
$ws.Range("N85").Value = -8817.6
$ws.Range("H100").Value = 3292.9412
$ws.Range("J100").Value = 4222
$ws.Range("L100").Value = 4222
$ws.Range("N100").Value = -5304
$ws.Range("H113").Value = 10361.0625
$ws.Range("I113").Value = 9290.615
$ws.Range("K113").Value = 9290.615
$ws.Range("M113").Value = -7120.615
$ws.Range("H122").Value = 7485.5557
$ws.Range("I122").Value = 8494
$ws.Range("J122").Value = 6225
$ws.Range("K122").Value = 25482
$ws.Range("L122").Value = 18675
$ws.Range("M122").Value = -23032
$ws.Range("N122").Value = -23575
$ws.Range("H125").Value = 78997.5
$ws.Range("J125").Value = 78997.5
$ws.Range("L125").Value = 78997.5
$ws.Range("N125").Value = -88837.5
$ws.Range("H136").Value = 2647.9062
$ws.Range("I136").Value = 1864.619
$ws.Range("J136").Value = 4143.273
$ws.Range("K136").Value = 5593.857
$ws.Range("L136").Value = 12429.819
$ws.Range("M136").Value = -3043.857
$ws.Range("N136").Value = -17529.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3500
$ws.Range("I14").Value = 3500
$ws.Range("K14").Value = 3500
$ws.Range("M14").Value = -3332
$ws.Range("H41").Value = 12697
$ws.Range("I41").Value = 16594.334
$ws.Range("J41").Value = 9774
$ws.Range("K41").Value = 16594.334
$ws.Range("L41").Value = 9774
$ws.Range("M41").Value = -16204.334
$ws.Range("N41").Value = -10554
$ws.Range("H96").Value = 2669554
$ws.Range("I96").Value = 6224847.5
$ws.Range("J96").Value = 3083.75
$ws.Range("K96").Value = 6224847.5
$ws.Range("L96").Value = 3083.75
$ws.Range("M96").Value = -6223474.5
$ws.Range("N96").Value = -5829.75
$ws.Range("H122").Value = 4580.8823
$ws.Range("I122").Value = 4537.879
$ws.Range("K122").Value = 13613.637
$ws.Range("M122").Value = -11163.637
$ws.Range("H123").Value = 39125
$ws.Range("J123").Value = 39125
$ws.Range("L123").Value = 39125
$ws.Range("N123").Value = -48925
$ws.Range("H136").Value = 4074.4
$ws.Range("I136").Value = 1217.5
$ws.Range("J136").Value = 15502
$ws.Range("K136").Value = 3652.5
$ws.Range("L136").Value = 46506
$ws.Range("M136").Value = -1102.5
$ws.Range("N136").Value = -51606
